$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH120"
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.ThemeColor = 1

$ws.Range("C2:H2").Font.Name = "Calibri"
$ws.Range("C2:H2").Font.ThemeColor = 1

$ws.Range("C2").Value = "MUSLIMS AGAINST APARTHEID CONFERENCE, IMAM W. DEEN MOHAMMED, AMERICAN MUSLIM COUNCIL, APARTHEID LAWS & INSTITUTIONS STILL IN EFFECT"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

$ws.Range("A5").Select()
